$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ row=2; E=3; F=1; G=2.802375666666666; H=8.407126999999999; I=0.004883830317835578; J=0.004883830317835578; M=0.649981; N=1.949943; O=0.04902821847542373; P=0.04902821847542374; Q=1.821490938195667; R=16.393418443761; S=0.0002394454998197408; T=0.0002394454998197409 },
    @{ row=3; E=3; F=1; G=2.802375666666666; H=8.407126999999999; I=0.004883830317835578; J=0.004883830317835578; O=0.1453376157926368; P=0.1453376157926368; Q=5.399566991771333; R=48.59610292594199; S=0.0007098042543300187; T=0.0007098042543300187 },
    @{ row=4; E=3; F=1; G=2.802375666666666; H=8.407126999999999; I=0.004883830317835578; J=0.004883830317835578; M=7.907236000000001; N=23.721708; O=0.596444656297239; P=0.5964446562972391; Q=22.15904575699067; R=199.431411812916; S=0.002912934495335476; T=0.002912934495335477 },
    @{ row=5; E=3; F=1; G=2.802375666666666; H=8.407126999999999; I=0.004883830317835578; J=0.004883830317835578; M=2.773284666666667; N=8.319853999999999; O=0.2091895094347004; P=0.2091895094347004; Q=7.771785466606444; R=69.94606919945799; S=0.001021646068350341; T=0.001021646068350342 },
    @{ row=6; G=543.3469646666667; I=0.9469160079809679; J=0.946916007980968; M=0.649981; N=1.949943; O=0.04902821847542373; P=0.04902821847542374; Q=353.1652034410047; R=3178.486830969042; S=0.04642560491716698; T=0.04642560491716699 },
    @{ row=7; G=543.3469646666667; I=0.9469160079809679; J=0.946916007980968; O=0.1453376157926368; P=0.1453376157926368; Q=1046.91115127437; S=0.1376225149558354; T=0.1376225149558354 },
    @{ row=8; G=543.3469646666667; I=0.9469160079809679; J=0.946916007980968; M=7.907236000000001; N=23.721708; O=0.596444656297239; P=0.5964446562972391; Q=4296.372679502996; R=38667.35411552696; S=0.564782992922562; T=0.5647829929225622 },
    @{ row=9; G=543.3469646666667; I=0.9469160079809679; J=0.946916007980968; M=2.773284666666667; N=8.319853999999999; O=0.2091895094347004; P=0.2091895094347004; Q=1506.855805789942; R=13561.70225210948; S=0.1980848951854035; T=0.1980848951854036 },
    @{ row=10; G=25.919625; H=77.758875; I=0.04517133513098909; J=0.0451713351309891; M=0.649981; N=1.949943; O=0.04902821847542373; P=0.04902821847542374; Q=16.847263777125; R=151.625373994125; S=0.002214670087628717; T=0.002214670087628717 },
    @{ row=11; G=25.919625; H=77.758875; I=0.04517133513098909; J=0.0451713351309891; O=0.1453376157926368; P=0.1453376157926368; Q=49.94146689675; R=449.47320207075; S=0.006565094150108132; T=0.006565094150108133 },
    @{ row=12; G=25.919625; H=77.758875; I=0.04517133513098909; J=0.0451713351309891; M=7.907236000000001; N=23.721708; O=0.596444656297239; P=0.5964446562972391; Q=204.9525919065; R=1844.5733271585; S=0.02694220145669019; T=0.0269422014566902 },
    @{ row=13; G=25.919625; H=77.758875; I=0.04517133513098909; J=0.0451713351309891; M=2.773284666666667; N=8.319853999999999; O=0.2091895094347004; P=0.2091895094347004; Q=71.88249857824999; R=646.94248720425; S=0.009449369436562056; T=0.00944936943656206 },
    @{ row=14; G=1.737961666666666; H=5.213884999999999; I=0.003028826570207414; J=0.003028826570207415; M=0.649981; N=1.949943; O=0.04902821847542373; P=0.04902821847542374; Q=1.129642062061667; R=10.166778558555; S=0.0001484979708082975; T=0.0001484979708082975 },
    @{ row=15; G=1.737961666666666; H=5.213884999999999; I=0.003028826570207414; J=0.003028826570207415; O=0.1453376157926368; P=0.1453376157926368; Q=3.348673256023333; R=30.13805930421; S=0.0004402024323633352; T=0.0004402024323633353 },
    @{ row=16; G=1.737961666666666; H=5.213884999999999; I=0.003028826570207414; J=0.003028826570207415; M=7.907236000000001; N=23.721708; O=0.596444656297239; P=0.5964446562972391; Q=13.74247305728667; R=123.68225751558; S=0.001806527422651306; T=0.001806527422651307 },
    @{ row=17; G=1.737961666666666; H=5.213884999999999; I=0.003028826570207414; J=0.003028826570207415; M=2.773284666666667; N=8.319853999999999; O=0.2091895094347004; P=0.2091895094347004; Q=4.81986244142111; R=43.37876197279; S=0.0006335987443844751; T=0.0006335987443844754 }
)

$colIndex = @{
    'A' = 1
    'B' = 2
    'C' = 3
    'D' = 4
    'E' = 5
    'F' = 6
    'G' = 7
    'H' = 8
    'I' = 9
    'J' = 10
    'K' = 11
    'L' = 12
    'M' = 13
    'N' = 14
    'O' = 15
    'P' = 16
    'Q' = 17
    'R' = 18
    'S' = 19
    'T' = 20
}

foreach ($entry in $updates) {
    $r = $entry.row
    foreach ($key in $entry.Keys) {
        if ($key -eq "row") { continue }
        $c = $colIndex[$key]
        $ws.Cells.Item($r, $c).Value = $entry[$key]
    }
}

Write-Host "Updated NATMI TPM values"
